# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted into the data table at row 245,
# pushing the existing rows 245-303 down to 246-304.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 245 (shifts 245..303 -> 246..304).
$ws.Rows.Item(245).Insert()

# Populate the newly inserted row 245 with the new record's data.
$ws.Range("A245").Value = 5
$ws.Range("B245").Value = "Macroferia Regional de Talca"
$ws.Range("C245").Value = "Maule"
$ws.Range("D245").Value = 44463
$ws.Range("E245").Value = 7
$ws.Range("F245").Value = "Fruta"
$ws.Range("G245").Value = 100104
$ws.Range("H245").Value = "Frutos de pepita"
$ws.Range("I245").Value = 100104005
$ws.Range("J245").Value = "Pera"
$ws.Range("K245").Value = "Packham's Triumph"
$ws.Range("L245").Value = "Especial"
$ws.Range("M245").Value = 300
$ws.Range("N245").Value = 12000
$ws.Range("O245").Value = 12000
$ws.Range("P245").Value = 12000
$ws.Range("Q245").Value = "`$/bandeja 18 kilos granel"
$ws.Range("R245").Value = "Provincia de Linares"
$ws.Range("S245").Value = 667
$ws.Range("T245").Value = 18
